$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices formatted as plain text (e.g. "226.89" or
# thousand-grouped "25.804.05"). Excel auto-detects simple decimal-looking
# strings as numbers, so we force those specific cells to Text format first
# to preserve the original string formatting (trailing zeros, etc.).
$ws.Range("D2").Value = "25.804.05"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.736.48"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.89"
$ws.Range("E5").Value = "  -4.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5127"
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2676"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.19"
$ws.Range("E9").Value = "  -5.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06071"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").Value = "1.735.66"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06994"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.08"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6233"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.485"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.17"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "25.820.28"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.39"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006508"
$ws.Range("E21").Value = "  -4.82%  "
$ws.Range("D22").Value = "1.955.61"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.023"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.346"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.060"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "136.42"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.501"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.91"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.43"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08280"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.600"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.345"
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04400"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.606"
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9681"
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5929"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01561"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.908"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9989"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3773"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7259"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.821"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05479"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.208"
$ws.Range("E47").Value = "  +4.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1093"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.59"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.56"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("E51").Value = "  +0.07%  "
